# Generate Report for Handoff
#
# The localization status workbook has three sheets:
#   - "Overview" : summary sheet (columns A-G)
#   - "zh-cn"    : per-file status table for the zh-cn locale
#   - "de-de"    : per-file status table for the de-de locale
#
# The handoff run moves the single tracked file from "In Translation" to
# "Ready for handoff" and refreshes the associated timestamps:
#   - Overview!G2          (Latest HO Xliff Generate Date)   -> 2016-10-14 07:06:29
#   - zh-cn!H2              (Latest Handoff Datetime)         -> 2016-10-14 07:06:18
#   - de-de!H2              (Latest Handoff Datetime)         -> 2016-10-14 07:06:29
#
# Updating the status text makes it longer ("In Translation" -> "Ready for
# handoff"), so the Status column on each sheet widens to fit the new text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-10-14 07:06:29"

# Columns E and F (zh-cn / de-de status) grow to fit the new text.
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333

# --- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-10-14 07:06:18"

# Status column (C) grows to fit the new text.
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333

# --- de-de sheet -----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-10-14 07:06:29"

# Status column (C) grows to fit the new text.
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
